# Optuna Attempt (go back with original)
# Update forecast metrics on "Forecast Comparison" and "Summary" sheets.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: Inventory Coverage (H) and Seasonality Index (L) ---

$wsForecast.Range("H2").Value = 2.5
$wsForecast.Range("L2").Value = 0.84

$wsForecast.Range("H3").Value = 1.5
$wsForecast.Range("L3").Value = 0.93

$wsForecast.Range("H4").Value = 0.5
$wsForecast.Range("L4").Value = 1.02

$wsForecast.Range("L5").Value = 0.99

$wsForecast.Range("L6").Value = 0.87

$wsForecast.Range("L7").Value = 0.98

$wsForecast.Range("L9").Value = 0.91

$wsForecast.Range("L10").Value = 1.02

$wsForecast.Range("L11").Value = 0.97

$wsForecast.Range("L12").Value = 1.17

$wsForecast.Range("L13").Value = 0.97

$wsForecast.Range("L14").Value = 1.02

$wsForecast.Range("L15").Value = 0.9

$wsForecast.Range("L16").Value = 0.83

$wsForecast.Range("L17").Value = 1.14

# --- Summary sheet: Total Forecast rows ---
# These cells hold numeric-looking text (inline/shared strings), so force
# Text formatting before assigning to keep them as strings rather than
# letting Excel auto-convert them to numbers.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "5"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "3"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "2"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "0"
